# Weekly data update: insert a new price record at the top of the Kiwi
# price-history table (row 556), pushing the existing historical rows
# (556-628) down by one (to 557-629).
#
# The workbook's single sheet stores one row per market observation,
# most-recent-first, for "Fruta, Terminal La Palmera de La Serena - Kiwi".
# Columns: A Mercado ID, B Mercado, C Región, D Fecha, E Codreg, F Tipo,
#          G Producto ID, H Producto, I Categoría ID, J Categoría,
#          K Variedad, L Calidad, M Volumen, N Precio mínimo,
#          O Precio máximo, P Precio promedio ponderado,
#          Q Unidad de comercialización, R Origen, S Precio $/Kg,
#          T Kg / unidad.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 556; this shifts the
# existing rows 556:628 down to 557:629 and extends the sheet dimension
# from A1:T628 to A1:T629 automatically.
$ws.Rows("556:556").Insert()

$newRow = 556

$ws.Cells.Item($newRow, 1).Value = 8
$ws.Cells.Item($newRow, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($newRow, 3).Value = "Coquimbo"
$ws.Cells.Item($newRow, 4).Value = 45142
$ws.Cells.Item($newRow, 5).Value = 4
$ws.Cells.Item($newRow, 6).Value = "Fruta"
$ws.Cells.Item($newRow, 7).Value = 100101
$ws.Cells.Item($newRow, 8).Value = "Berries"
$ws.Cells.Item($newRow, 9).Value = 100101007
$ws.Cells.Item($newRow, 10).Value = "Kiwi"
$ws.Cells.Item($newRow, 11).Value = "Hayward"
$ws.Cells.Item($newRow, 12).Value = "Primera"
$ws.Cells.Item($newRow, 13).Value = 10
$ws.Cells.Item($newRow, 14).Value = 310000
$ws.Cells.Item($newRow, 15).Value = 320000
$ws.Cells.Item($newRow, 16).Value = 315000
$ws.Cells.Item($newRow, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item($newRow, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($newRow, 19).Value = 700
$ws.Cells.Item($newRow, 20).Value = 450
